$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Benutzer" table -> renamed/reworked into "User" table ---
$ws.Range("A1").Value = "User"
$ws.Range("B2").Value = "userID (Varchar)"
$ws.Range("C2").Value = "password (Varchar)"
$ws.Range("D2").Value = "isFarmer (int)"

# --- "Offers" table gains two more columns (distance, date) and the
#     second column now references the renamed user table ---
$ws.Range("C7").Value = "userID(VARCHAR) REFERENCES user(userID)"
$ws.Range("D7").Value = "distance(int)"
$ws.Range("E7").Value = "date(varchar)"

# --- column widths: widen column C, and size the new column E ---
$ws.Columns.Item(3).ColumnWidth = 35
$ws.Columns.Item(5).ColumnWidth = 12.6666667

# --- selection moves from D18 to E18 to track the new rightmost column ---
$ws.Range("E18").Select()
